# Apply the "TMML-Changes-Suggested" update:
#  - Rename header D1 "Impact" -> "Impact on design"
#  - Add two new trailing columns: "Status" (F) and "Man Hours" (G)
#  - Fix Change Type on row 18 (D18): "Medium" -> "Small"
#  - Add two new data rows (22, 23) to the change log
#  - Extend Table1 to cover the new columns/rows
#  - Update sheet view (scroll/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Header row updates
# ---------------------------------------------------------------
$ws.Range("D1").Value = "Impact on design"

# ---------------------------------------------------------------
# 2. Fix existing data: row 18 "Change Type" Medium -> Small
# ---------------------------------------------------------------
$ws.Range("D18").Value = "Small"

# ---------------------------------------------------------------
# 3. New rows 22 and 23
# ---------------------------------------------------------------
$ws.Range("A22").Value = "Machine Master"
$ws.Range("B22").Value = "Part process sequnce can have machine or machine group."
$ws.Range("C22").Value = "small"
$ws.Range("D22").Value = "Minor"
$ws.Range("E22").Value = ""

$ws.Range("A23").Value = "part Number "
$ws.Range("B23").Value = "Once finished, there should be another confirmation apI call where user scan part number and sends the part number and its quantity"
$ws.Range("C23").Value = "Big"
$ws.Range("D23").Value = "Big"
$ws.Range("E23").Value = "Was not part of earlier scope"

# Row heights matching the rest of the table (17 = single line, 34 = two lines)
$ws.Rows.Item(22).RowHeight = 17
$ws.Rows.Item(23).RowHeight = 34

# Formatting to match the surrounding (non-header) table rows: 12pt Calibri,
# left/center aligned; column B wraps text.
$row22 = $ws.Range("A22:E22")
$row22.Font.Size = 12
$row22.Font.Bold = $false
$row22.HorizontalAlignment = -4131
$row22.VerticalAlignment = -4108
$row22.WrapText = $false
$ws.Range("B22").WrapText = $true

$row23 = $ws.Range("A23:E23")
$row23.Font.Size = 12
$row23.Font.Bold = $false
$row23.HorizontalAlignment = -4131
$row23.VerticalAlignment = -4108
$row23.WrapText = $false
$ws.Range("B23").Font.Bold = $true
$ws.Range("B23").WrapText = $true

# ---------------------------------------------------------------
# 4. Add new "Status" and "Man Hours" table columns (F, G)
# ---------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null

$ws.Range("F1").Value = "Status"
$ws.Range("G1").Value = "Man Hours"

$header = $ws.Range("F1:G1")
$header.Font.Size = 12
$header.Font.Bold = $true
$header.HorizontalAlignment = -4131
$header.VerticalAlignment = -4108
$header.WrapText = $false

# Match formatting of the rest of column F/G (regular, non-wrapped 12pt cells)
$fg = $ws.Range("F2:G23")
$fg.Font.Size = 12
$fg.Font.Bold = $false
$fg.HorizontalAlignment = -4131
$fg.VerticalAlignment = -4108
$fg.WrapText = $false

# Grow the table down to the new data rows (22 and 23)
$lo.Resize($ws.Range("A1:G23")) | Out-Null

# ---------------------------------------------------------------
# 5. Sheet view: scroll so column B is the left-most visible column,
#    and select G2 (matches the recorded view state of the edit)
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G2").Select() | Out-Null
